$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 25 (Designator "R1") -- unnecessary resistor removed.
$ws.Rows(25).Delete()

# Move the active selection to B27 (per the saved view state in the edit).
$ws.Range("B27").Select()
